$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1001
$ws.Range("J43").Value = 1001
$ws.Range("L43").Value = 1001
$ws.Range("N43").Value = -1139
$ws.Range("H62").Value = 250008460
$ws.Range("I62").Value = 500001440
$ws.Range("J62").Value = 15490
$ws.Range("K62").Value = 500001440
$ws.Range("L62").Value = 15490
$ws.Range("M62").Value = -500000816
$ws.Range("N62").Value = -16738
$ws.Range("H65").Value = 250008460
$ws.Range("I65").Value = 500001440
$ws.Range("J65").Value = 15490
$ws.Range("K65").Value = 2500007200
$ws.Range("L65").Value = 77450
$ws.Range("M65").Value = -2500004080
$ws.Range("N65").Value = -83690
$ws.Range("H100").Value = 2412.7
$ws.Range("I100").Value = 1432.8334
$ws.Range("K100").Value = 1432.8334
$ws.Range("M100").Value = -891.8334
$ws.Range("H133").Value = 59440
$ws.Range("J133").Value = 59440
$ws.Range("L133").Value = 59440
$ws.Range("N133").Value = -69560
$ws.Range("H137").Value = 2021.591
$ws.Range("I137").Value = 1848.1333
$ws.Range("J137").Value = 2393.2856
$ws.Range("K137").Value = 5544.3999
$ws.Range("L137").Value = 7179.8568
$ws.Range("M137").Value = -2994.3999
$ws.Range("N137").Value = -12279.8568
$ws.Range("H138").Value = 1603.3334
$ws.Range("J138").Value = 2102.5676
$ws.Range("L138").Value = 6307.702799999999
$ws.Range("N138").Value = -16587.7028

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 327923.2
$ws.Range("I2").Value = 464101.66
$ws.Range("J2").Value = 1094.8
$ws.Range("K2").Value = 464101.66
$ws.Range("L2").Value = 1094.8
$ws.Range("M2").Value = -463988.66
$ws.Range("N2").Value = -1320.8
$ws.Range("H32").Value = 4789
$ws.Range("I32").Value = 3082.8794
$ws.Range("J32").Value = 17158.375
$ws.Range("K32").Value = 3082.8794
$ws.Range("L32").Value = 17158.375
$ws.Range("M32").Value = -2795.8794
$ws.Range("N32").Value = -17732.375
$ws.Range("H61").Value = 6360.2173
$ws.Range("I61").Value = 7438.9287
$ws.Range("K61").Value = 7438.9287
$ws.Range("M61").Value = -7226.9287
$ws.Range("H110").Value = 250.11111
$ws.Range("I110").Value = 250.11111
$ws.Range("K110").Value = 250.11111
$ws.Range("M110").Value = 1794.88889
$ws.Range("H116").Value = 327923.2
$ws.Range("I116").Value = 464101.66
$ws.Range("J116").Value = 1094.8
$ws.Range("K116").Value = 464101.66
$ws.Range("L116").Value = 1094.8
$ws.Range("M116").Value = -461807.66
$ws.Range("N116").Value = -5682.8
$ws.Range("H136").Value = 6360.2173
$ws.Range("I136").Value = 7438.9287
$ws.Range("K136").Value = 22316.7861
$ws.Range("M136").Value = -19766.7861

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 327923.2
$ws.Range("I3").Value = 464101.66
$ws.Range("J3").Value = 1094.8
$ws.Range("K3").Value = 464101.66
$ws.Range("L3").Value = 1094.8
$ws.Range("M3").Value = -463987.66
$ws.Range("N3").Value = -1322.8
$ws.Range("H134").Value = 5131.3667
$ws.Range("I134").Value = 5844.478
$ws.Range("K134").Value = 17533.434
$ws.Range("M134").Value = -14998.434

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2297.9
$ws.Range("I31").Value = 2725
$ws.Range("J31").Value = 2250.4443
$ws.Range("K31").Value = 2725
$ws.Range("L31").Value = 2250.4443
$ws.Range("M31").Value = -2430
$ws.Range("N31").Value = -2840.4443
$ws.Range("H34").Value = 2297.9
$ws.Range("I34").Value = 2725
$ws.Range("J34").Value = 2250.4443
$ws.Range("K34").Value = 2725
$ws.Range("L34").Value = 2250.4443
$ws.Range("M34").Value = -2523
$ws.Range("N34").Value = -2654.4443
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = $null
$ws.Range("M80").Value = $null
$ws.Range("N80").Value = 0
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = $null
$ws.Range("M83").Value = $null
$ws.Range("N83").Value = 0
$ws.Range("H132").Value = 3060.4666
$ws.Range("I132").Value = 1737.375
$ws.Range("K132").Value = 5212.125
$ws.Range("M132").Value = -2682.125
$ws.Range("H134").Value = 3783.2
$ws.Range("I134").Value = 3429.4285
$ws.Range("K134").Value = 10288.2855
$ws.Range("M134").Value = -7753.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 10702.857
$ws.Range("I120").Value = 8930
$ws.Range("J120").Value = 10998.333
$ws.Range("K120").Value = 26790
$ws.Range("L120").Value = 32994.999
$ws.Range("M120").Value = -21952
$ws.Range("N120").Value = -42670.999
$ws.Range("H129").Value = 91339.625
$ws.Range("J129").Value = 145597.8
$ws.Range("L129").Value = 436793.4
$ws.Range("N129").Value = -446793.4
$ws.Range("H131").Value = 13701.148
$ws.Range("J131").Value = 14198.692
$ws.Range("L131").Value = 42596.076
$ws.Range("N131").Value = -52676.076
$ws.Range("H136").Value = 1551
$ws.Range("I136").Value = 1551
$ws.Range("K136").Value = 4653
$ws.Range("M136").Value = 447

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 631.35
$ws.Range("I97").Value = 625.1053000000001
$ws.Range("J97").Value = 750
$ws.Range("K97").Value = 625.1053000000001
$ws.Range("L97").Value = 750
$ws.Range("M97").Value = -129.1053000000001
$ws.Range("N97").Value = -1742
$ws.Range("H113").Value = 1355
$ws.Range("I113").Value = 1132
$ws.Range("J113").Value = 1503.6666
$ws.Range("K113").Value = 1132
$ws.Range("L113").Value = 1503.6666
$ws.Range("M113").Value = 1038
$ws.Range("N113").Value = -5843.6666
$ws.Range("H132").Value = 3500555.2
$ws.Range("I132").Value = 6412852
$ws.Range("J132").Value = 5799.2
$ws.Range("K132").Value = 19238556
$ws.Range("L132").Value = 17397.6
$ws.Range("M132").Value = -19236026
$ws.Range("N132").Value = -22457.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 10865.143
$ws.Range("J43").Value = 10865.143
$ws.Range("L43").Value = 10865.143
$ws.Range("N43").Value = -11251.143
$ws.Range("H55").Value = 510.89474
$ws.Range("J55").Value = 491.55554
$ws.Range("L55").Value = 491.55554
$ws.Range("N55").Value = -837.5555400000001
$ws.Range("H122").Value = 6383.278
$ws.Range("I122").Value = 6263.091
$ws.Range("K122").Value = 18789.273
$ws.Range("M122").Value = -16339.273
$ws.Range("H136").Value = 8000
$ws.Range("I136").Value = 6000
$ws.Range("K136").Value = 18000
$ws.Range("M136").Value = -15450

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = $null
$ws.Range("N31").Value = 0
$ws.Range("H70").Value = 29110
$ws.Range("J70").Value = 29110
$ws.Range("L70").Value = 29110
$ws.Range("N70").Value = -29740
$ws.Range("H73").Value = 29110
$ws.Range("J73").Value = 29110
$ws.Range("L73").Value = 29110
$ws.Range("N73").Value = -31294
$ws.Range("H96").Value = 1298.25
$ws.Range("I96").Value = 901.5
$ws.Range("J96").Value = 1430.5
$ws.Range("K96").Value = 901.5
$ws.Range("L96").Value = 1430.5
$ws.Range("M96").Value = 471.5
$ws.Range("N96").Value = -4176.5
$ws.Range("H122").Value = 42180.473
$ws.Range("I122").Value = 52842.4
$ws.Range("K122").Value = 158527.2
$ws.Range("M122").Value = -156077.2
$ws.Range("H126").Value = 1966.8462
$ws.Range("I126").Value = 1972
$ws.Range("J126").Value = 1958.6
$ws.Range("K126").Value = 5916
$ws.Range("L126").Value = 5875.799999999999
$ws.Range("M126").Value = -3446
$ws.Range("N126").Value = -10815.8
$ws.Range("H132").Value = 1072.2759
$ws.Range("I132").Value = 824.61365
$ws.Range("K132").Value = 2473.84095
$ws.Range("M132").Value = 56.15905000000021
